# Logins.xlsx - test version of the login screen.
# Remove the two unused placeholder sheets and populate Plan1 with the
# login table (user / password / email / name / age), then leave the
# selection where the author left it (E10, below the data).

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Drop the empty extra sheets - only Plan1 survives.
[void]$wb.Worksheets.Item("Plan2").Delete()
[void]$wb.Worksheets.Item("Plan3").Delete()

# Fill in rows 2-6 first (the body of the table)...
$ws.Range("A2").Value = "Azimba"
$ws.Range("B2").Value = "xabalabaluba"
$ws.Range("C2").Value = "FCV@brotherbund.com"
$ws.Range("D2").Value = "Fabio"
$ws.Range("E2").Value = 25

$ws.Range("A3").Value = "Goloko"
$ws.Range("B3").Value = "zoiao69"
$ws.Range("C3").Value = "ASS@brotherbund.com"
$ws.Range("D3").Value = "Danilo"
$ws.Range("E3").Value = 19

$ws.Range("A4").Value = "Baiuca1"
$ws.Range("B4").Value = "zonasul99"
$ws.Range("C4").Value = "MOE@brotherbund.com"
$ws.Range("D4").Value = "Moe"
$ws.Range("E4").Value = 39

$ws.Range("A5").Value = "jobijoba"
$ws.Range("B5").Value = "bamboleo"
$ws.Range("C5").Value = "GSK@brotherbund.com"
$ws.Range("D5").Value = "Carlos"
$ws.Range("E5").Value = 20

$ws.Range("A6").Value = "user06"
$ws.Range("B6").Value = "variant01"
$ws.Range("C6").Value = "FBX@brotherbund.com"
$ws.Range("D6").Value = "Aline"
$ws.Range("E6").Value = 23

# ...then add row 1 last (matches the author's edit order / shared-string layout).
$ws.Range("A1").Value = "user01"
$ws.Range("B1").Value = "variant01"
$ws.Range("C1").Value = "FBZ@brotherbund.com"
$ws.Range("D1").Value = "Ana"
$ws.Range("E1").Value = 18

# Leave the cursor where the author left it when committing.
[void]$ws.Range("E10").Select()
